# Add the Shoe Size (D) and Eye Color (E) data that was missing from the
# "Data" sheet. Headers already exist in D1/E1; only the data rows
# (2-15) need to be filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$shoeSize = @(8, 7, 6, 7.5, 9, 5.5, 9, 11, 5.5, 7, 10, 4.5, 5.5, 5)
$eyeColor = @("L", "G", "G", "R", "G", "G", "R", "R", "R", "L", "L", "R", "G", "L")

for ($i = 0; $i -lt $shoeSize.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $shoeSize[$i]
    $ws.Cells.Item($row, 5).Value = $eyeColor[$i]
}

$ws.Range("E16").Select()
